# Insert a new weekly record at row 229 (Hortaliza, Femacal de La Calera - Ciboulette).
# This pushes the existing rows 229-328 down to 230-329 and the sheet
# dimension grows from A1:R328 to A1:R329.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 229; Excel shifts rows 229:328
# down to 230:329 and the new row 229 inherits formatting from row 228.
$ws.Rows("229:229").Insert()

# Populate the new row 229 with the same constant attributes used by every
# other row in this data block, plus the new record's own Fecha/Volumen.
$ws.Range("A229").Value = 3
$ws.Range("B229").Value = "Femacal de La Calera"
$ws.Range("C229").Value = "Coquimbo"
$ws.Range("D229").Value = 44726
$ws.Range("E229").Value = 5
$ws.Range("F229").Value = 100112039
$ws.Range("G229").Value = "Ciboulette"
$ws.Range("H229").Value = "Sin especificar"
$ws.Range("I229").Value = "Primera"
$ws.Range("J229").Value = 160
$ws.Range("K229").Value = 1500
$ws.Range("L229").Value = 1500
$ws.Range("M229").Value = 1500
$ws.Range("N229").Value = "$/docena de atados"
$ws.Range("O229").Value = "Provincia de Quillota"
$ws.Range("P229").Value = 500
$ws.Range("Q229").Value = 3
$ws.Range("R229").Value = "Hortaliza"
